$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "595.79") must be
# protected with a temporary text format, otherwise Excel auto-converts the
# assigned string into a numeric value instead of keeping it as text.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "65.520.31"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "2.640.96"
$ws.Range("E3").Value = "  -1.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "595.79"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").Value = "155.68"
$ws.Range("E6").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  +3.15%  "

# Row 9
$ws.Range("E9").Value = "  +4.26%  "

# Row 10
$ws.Range("D10").Value = "0.396"
$ws.Range("E10").Value = "  +0.05%  "

# Row 11
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -2.40%  "

# Row 12
$ws.Range("E12").Value = "  +0.64%  "

# Row 13
$ws.Range("D13").Value = "28.69"
$ws.Range("E13").Value = "  -2.16%  "

# Row 14
$ws.Range("D14").Value = "0.0000195"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("D15").Value = "3.114.98"
$ws.Range("E15").Value = "  -1.10%  "

# Row 16
$ws.Range("D16").Value = "65.389.48"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "2.634.20"
$ws.Range("E17").Value = "  -1.41%  "

# Row 18
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  +0.21%  "

# Row 19
$ws.Range("D19").Value = "4.71"
$ws.Range("E19").Value = "  -1.87%  "

# Row 20
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  -1.15%  "

# Row 21
$ws.Range("D21").Value = "348.00"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("D23").Value = "68.80"
$ws.Range("E23").Value = "  -2.11%  "

# Row 24
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +2.67%  "

# Row 25
$ws.Range("D25").Value = "9.59"
$ws.Range("E25").Value = "  -2.03%  "

# Row 26
$ws.Range("D26").Value = "1.64"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").Value = "1.58"
$ws.Range("E27").Value = "  -1.90%  "

# Row 28
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -2.91%  "

# Row 29
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.56%  "

# Row 30
$ws.Range("D30").Value = "537.44"
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.14"
$ws.Range("E31").Value = "  -0.20%  "

# Row 32
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "7.83"
$ws.Range("E32").Value = "  -3.32%  "

# Row 33
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  -0.35%  "

# Row 34
$ws.Range("D34").Value = "6.40"
$ws.Range("E34").Value = "  -1.72%  "

# Row 35
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").Value = "  +0.78%  "

# Row 36
$ws.Range("D36").Value = "0.420"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37
$ws.Range("D37").Value = "20.30"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.01%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "152.97"
$ws.Range("E40").Value = "  -4.23%  "

# Row 42
$ws.Range("D42").Value = "159.77"
$ws.Range("E42").Value = "  -3.33%  "

# Row 43
$ws.Range("D43").Value = "4.06"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44
$ws.Range("E44").Value = "  +2.70%  "

# Row 45
$ws.Range("D45").Value = "0.0604"
$ws.Range("E45").Value = "  -1.50%  "

# Row 46
$ws.Range("D46").Value = "22.47"
$ws.Range("E46").Value = "  -2.14%  "

# Row 47
$ws.Range("E47").Value = "  -2.43%  "

# Row 48
$ws.Range("D48").Value = "0.0254"
$ws.Range("E48").Value = "  -2.40%  "

# Row 49
$ws.Range("D49").Value = "0.0990"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0250"
$ws.Range("E50").Value = "  +8.90%  "

# Row 51
$ws.Range("D51").Value = "19.56"
$ws.Range("E51").Value = "  -2.30%  "

# Restore default (Normal) style on the protected cells now that their
# text values are locked in, so no spurious formatting is introduced.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
